# Remove the obsolete "JSLHISAR" lookup row (row 99) from the instruments
# lookup table. Deleting the entire row shifts all subsequent rows up by
# one, removes the now-unused "JSLHISAR" shared string, and shrinks the
# used range from A1:B191 to A1:B190.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(99).Delete()

# Update the selected/active cell to match the author's final view state.
$ws.Range("G95").Select()
